# Added Columns to Asset List
# Adds three new columns (Contato, Analista, CSP) to the asset list sheet,
# filling in the "John" / "Mike" / "Anna" values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# Headers
$ws.Range("I1").Value = "Contato"
$ws.Range("J1").Value = "Analista"
$ws.Range("K1").Value = "CSP"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "John"
    $ws.Cells.Item($r, 10).Value = "Mike"
    $ws.Cells.Item($r, 11).Value = "Anna"
}

$ws.Range("N8").Select()
